# Updated legacy GSC export data.
#
# The oldest day in the "Chart" series (2025-08-24), which only ever held
# placeholder/blank counts, is dropped. Every later date's row shifts up by
# one so the series now runs 2025-08-25 .. 2025-11-20 (89 data rows instead
# of 90). The former last-row placeholder blank Impressions cell becomes an
# explicit numeric 0, matching how the other numeric columns are stored.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Remove the 2025-08-24 row; Excel shifts rows 3:90 up into rows 2:89.
$ws.Rows(2).Delete()

# The series now ends at row 89 (2025-11-20); make sure its Impressions
# value is a real number rather than the old blank placeholder.
$ws.Range("D89").Value = 0
